# Insert a new weekly price row at row 61 (Terminal La Palmera de La Serena - Ajo),
# shifting all existing rows 61:161 down to 62:162, and populate the new row 61
# with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("61:61").Insert()

$ws.Range("A61").Value = 8
$ws.Range("B61").Value = "Terminal La Palmera de La Serena"
$ws.Range("C61").Value = "Coquimbo"
$ws.Range("D61").Value = 44495
$ws.Range("E61").Value = 4
$ws.Range("F61").Value = 100112003
$ws.Range("G61").Value = "Ajo"
$ws.Range("H61").Value = "Chino"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 560
$ws.Range("K61").Value = 19000
$ws.Range("L61").Value = 20000
$ws.Range("M61").Value = 19500
$ws.Range("N61").Value = "$/caja 10 kilos"
$ws.Range("O61").Value = "China"
$ws.Range("P61").Value = 1950
$ws.Range("Q61").Value = 10
$ws.Range("R61").Value = "Hortaliza"
